# Apply the "SELF_INDUCTANCE_MODE" input-variable addition to the
# CONDUCTOR_operation sheet of template_conductor_definition.xlsx,
# and update the existing INDUCTANCE_MODE description / value.

$wb = $excel.ActiveWorkbook
$wsOperation = $wb.Worksheets.Item("CONDUCTOR_operation")

# --- Update existing INDUCTANCE_MODE row (row 8) ---------------------------
# New description: now 0 = analytical, 1 = approximated (default 1), so also
# bump the stored default value for the flag itself from 0 to 1.
$wsOperation.Range("D8").Value = "flag to select the method to evaluate the inductance. Possible values: 0 = analytical; 1 = approximated. Defaults to 1."
$wsOperation.Range("E8").Value = 1

# --- Insert the new SELF_INDUCTANCE_MODE row (row 9) ------------------------
# Clone the formatting of row 8 (font/fill/alignment/number format) onto the
# new row 9, then overwrite the actual contents.
$wsOperation.Range("A8:E8").Copy()
$wsOperation.Range("A9:E9").PasteSpecial(-4122)  # xlPasteFormats
$wsOperation.Range("A8:E8").Copy()
$wsOperation.Range("A9:E9").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$wsOperation.Rows.Item(9).RowHeight = 159.5

$wsOperation.Range("A9").Value = "SELF_INDUCTANCE_MODE"
$wsOperation.Range("B9").Value = "-"
$wsOperation.Range("C9").Value = "integer"
$wsOperation.Range("D9").Value = "flag to select the method to evaluate the self inductance. Possible values: 1 = mode 1; 2 = mode 2. Used only if flag INDUCTANCE_MODE is set to 1. Defaults to 2"
$wsOperation.Range("E9").Value = 2

# --- View-state bookkeeping: CONDUCTOR_operation becomes the active sheet --
$wsOperation.Activate()
$wsOperation.Range("G9").Select()
$excel.ActiveWindow.ScrollRow = 8
